$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 15.24 = 64558.08 pesos`n✅ 64558.08 pesos = 15.13 = 966.56 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update N10, O10 and O12 values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 65.59999999999999
$ws2.Range("O10").Value = 4235.01
$ws2.Range("O12").Value = 63.9
